$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.897.55'
$ws.Range('E2').Value = '  +7.10%  '
$ws.Range('D3').Value = '2.309.09'
$ws.Range('E3').Value = '  +5.76%  '
$ws.Range('E4').Value = '  -0.63%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '298.87'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '98.02'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +12.41%  '
$ws.Range('E7').Value = '  +2.53%  '
$ws.Range('E8').Value = '  -0.43%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.528'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +11.19%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '35.69'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +4.94%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '7.33'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +9.10%  '
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('D14').Value = '2.661.51'
$ws.Range('E14').Value = '  +5.57%  '
$ws.Range('D15').Value = '2.308.64'
$ws.Range('E15').Value = '  +1.82%  '
$ws.Range('E16').Value = '  +8.48%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.817'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +7.85%  '
$ws.Range('D18').Value = '46.784.15'
$ws.Range('E18').Value = '  +7.77%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.13'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +22.30%  '
$ws.Range('D20').Value = '0.0₃0940'
$ws.Range('E20').Value = '  +7.58%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.12'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +5.70%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '66.88'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +7.02%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '249.07'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +9.99%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.92'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +6.06%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '1.98'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +10.01%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '42.95'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +22.81%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.25'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +3.13%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.82'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +7.93%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '20.18'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +6.65%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '5.74'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +9.62%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '146.97'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.98%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0797'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +11.37%  '
$ws.Range('E34').Value = '  +6.26%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.11'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +8.38%  '
$ws.Range('E36').Value = '  +11.58%  '
$ws.Range('E37').Value = '  +3.87%  '
$ws.Range('E38').Value = '  +10.61%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '15.75'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +21.78%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.99'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +14.96%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '3.44'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +13.91%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0308'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +10.85%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.65%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.00'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +23.10%  '
$ws.Range('D45').Value = '1.836.52'
$ws.Range('E45').Value = '  +5.09%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '90.91'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +25.18%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.198'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +17.23%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '75.68'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +12.97%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '4.97'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +13.60%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '97.34'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +7.46%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '54.17'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +12.19%  '
